$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsOverview.Range("G2").Value = "2016-09-05 03:13:26"
$wsZhCn.Range("H2").Value = "2016-09-05 03:13:19"
$wsZhCn.Range("K2").Value = "2016-09-05 03:13:43"
$wsDeDe.Range("H2").Value = "2016-09-05 03:13:26"
$wsDeDe.Range("K2").Value = "2016-09-05 03:13:50"
